$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last three BOM rows (Balsa wood sheet/strip/stick items + their
# "Balsa USA" source) are being pulled from the bill of materials. Clear
# their contents but leave the row/cell formatting (wrap-text style on
# column A, hyperlink style on column F) intact, exactly like an
# Excel "Clear Contents" on A18:G20.
$ws.Range("A18:G20").ClearContents()

# Leave the selection where the user ended up after doing the edit.
$ws.Range("A20").Select()
